$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.974.15"
$ws.Range("E2").Value = "  -2.27%  "

$ws.Range("D3").Value = "2.097.09"
$ws.Range("E3").Value = "  -1.12%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  -0.85%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "346.44"
$ws.Range("E5").Value = "  +2.58%  "

$ws.Range("E6").Value = "  -0.82%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5152"
$ws.Range("E7").Value = "  -1.93%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4417"
$ws.Range("E8").Value = "  -3.08%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09272"
$ws.Range("E9").Value = "  +1.65%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "52.05"
$ws.Range("E10").Value = "  -5.14%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.168"
$ws.Range("E11").Value = "  -0.65%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "25.19"
$ws.Range("E12").Value = "  +2.83%  "

$ws.Range("D13").Value = "2.099.44"
$ws.Range("E13").Value = "  -1.22%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.731"
$ws.Range("E14").Value = "  -1.90%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.143"
$ws.Range("E15").Value = "  -0.19%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "99.59"
$ws.Range("E16").Value = "  +2.36%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001166"
$ws.Range("E17").Value = "  -0.87%  "

$ws.Range("E18").Value = "  -0.79%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "20.60"
$ws.Range("E19").Value = "  +5.71%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.06664"
$ws.Range("E20").Value = "  -0.49%  "

$ws.Range("E21").Value = "  -0.80%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.212"
$ws.Range("E22").Value = "  -1.78%  "

$ws.Range("D23").Value = "30.083.33"
$ws.Range("E23").Value = "  -2.18%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.60"
$ws.Range("E24").Value = "  -2.32%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.332"
$ws.Range("E25").Value = "  -1.17%  "

$ws.Range("D26").Value = "2.348.79"
$ws.Range("E26").Value = "  -1.08%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.94"
$ws.Range("E27").Value = "  -2.12%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.551"
$ws.Range("E28").Value = "  -0.54%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "162.25"
$ws.Range("E29").Value = "  -1.67%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.26"
$ws.Range("E30").Value = "  -1.14%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.167"
$ws.Range("E31").Value = "  -3.76%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1060"
$ws.Range("E32").Value = "  -1.42%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.636"
$ws.Range("E33").Value = "  -1.84%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.220"
$ws.Range("E34").Value = "  -2.57%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.956"
$ws.Range("E35").Value = "  +0.30%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.187"
$ws.Range("E36").Value = "  +5.06%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "10.06"
$ws.Range("E37").Value = "  -5.58%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02560"
$ws.Range("E38").Value = "  -3.07%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06773"
$ws.Range("E39").Value = "  -1.67%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2271"
$ws.Range("E40").Value = "  -2.38%  "

$ws.Range("B41").Value = "Aptos"
$ws.Range("C41").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.48"
$ws.Range("E41").Value = "  -1.67%  "

$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6895"
$ws.Range("E42").Value = "  -0.45%  "

$ws.Range("E43").Value = "  +3.95%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6620"
$ws.Range("E44").Value = "  +1.92%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.26"
$ws.Range("E45").Value = "  -6.16%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.276"
$ws.Range("E46").Value = "  -2.03%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.630"
$ws.Range("E47").Value = "  -1.83%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.00000000353"
$ws.Range("E48").Value = "  -5.20%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.218"
$ws.Range("E49").Value = "  -3.07%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "81.92"
$ws.Range("E50").Value = "  -1.90%  "

$ws.Range("E51").Value = "  -1.49%  "
